$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.690.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.89%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.562.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.01%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'565.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.83%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.43%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.563.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.113"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.28%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.76%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'TRON"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.156"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'Cardano"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.372"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.027.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.94%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.582.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.577.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.53%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.00%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.69%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'332.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'66.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.57%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +5.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000108"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'WrappedeETH"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.709.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.42%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Fetch.AI"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.86%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Bittensor"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'556.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.78%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Aptos"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.12%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'Kaspa"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'0.157"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.88%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'PancakeSwap"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'1.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'RenderToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'6.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'NEARProtocol"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'5.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'FirstDigitalUSD"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.11%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'PolygonEcosystemToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.75%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'EthereumClassic"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'19.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.98%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Monero"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'151.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.81%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Stacks"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'USDe"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'dogwifhat"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.18%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Aave"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'154.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'23.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Filecoin"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'3.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Hedera"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.0568"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Mantle"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.617"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.89%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Stellar"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'VeChain"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0242"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.63%  "
$ws.Range("E51").Style = "Normal"
